$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg'
$ws.Range("G3").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G4").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'
$ws.Range("G5").Value = 'lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Nadagouda@med.asu.edu.eg'
$ws.Range("G7").Value = 'Shimaa.ashraf@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg'
$ws.Range("G8").Value = 'Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G13").Value = 'marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg'
$ws.Range("G16").Value = 'aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G18").Value = 'Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg'
$ws.Range("G19").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G21").Value = 'lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Nadagouda@med.asu.edu.eg'
$ws.Range("G23").Value = 'Shimaa.ashraf@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg'
$ws.Range("G24").Value = 'Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G29").Value = 'marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg'
$ws.Range("G32").Value = 'aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G34").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G35").Value = 'Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G50").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G51").Value = 'Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G66").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G67").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G69").Value = 'NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg'
$ws.Range("G80").Value = 'ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
$ws.Range("G82").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G83").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range("G85").Value = 'NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg'
$ws.Range("G96").Value = 'ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
$ws.Range("G98").Value = 'Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg'
$ws.Range("G103").Value = 'Shimaa.ashraf@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg'
$ws.Range("G114").Value = 'Mohammedeltanany@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg'
$ws.Range("G116").Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'
$ws.Range("G119").Value = 'Shimaa.ashraf@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg'
